$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bulgaria First League")

    $ws.Range("B8").Value = 6627736
    $ws.Range("E8").Value = "Botev Plovdiv"
    $ws.Range("F8").Value = "Arda Kardzhali"
    $ws.Range("G8").Value = 0
    $ws.Range("H8").Value = 3
    $ws.Range("J8").Value = 1
    $ws.Range("K8").Value = "A"
    $ws.Range("L8").Value = 5.25
    $ws.Range("M8").Value = 3.6
    $ws.Range("N8").Value = 1.571
    $ws.Range("O8").Value = 26
    $ws.Range("P8").Value = 11
    $ws.Range("Q8").Value = 1.083
    $ws.Range("R8").Value = 2.5
    $ws.Range("S8").Value = 1.825
    $ws.Range("T8").Value = 2.025
    $ws.Range("U8").Value = 3.25
    $ws.Range("V8").Value = 2
    $ws.Range("W8").Value = 1.85
    $ws.Range("X8").Value = -1
    $ws.Range("Z8").Value = 0.08299999999999996
    $ws.Range("AA8").Value = -1
    $ws.Range("AB8").Value = 1.025
    $ws.Range("AD8").Value = 0.425
    $ws.Range("B9").Value = 6627737
    $ws.Range("E9").Value = "Slavia Sofia"
    $ws.Range("F9").Value = "Lokomotiv 1929 Sofia"
    $ws.Range("G9").Value = 2
    $ws.Range("H9").Value = 0
    $ws.Range("J9").Value = 0
    $ws.Range("K9").Value = "H"
    $ws.Range("L9").Value = 1.5
    $ws.Range("M9").Value = 3.75
    $ws.Range("N9").Value = 6.5
    $ws.Range("O9").Value = 1.444
    $ws.Range("P9").Value = 4.333
    $ws.Range("Q9").Value = 8
    $ws.Range("R9").Value = -1.25
    $ws.Range("S9").Value = 2
    $ws.Range("T9").Value = 1.85
    $ws.Range("U9").Value = 2.25
    $ws.Range("V9").Value = 1.875
    $ws.Range("W9").Value = 1.975
    $ws.Range("X9").Value = 0.444
    $ws.Range("Z9").Value = -1
    $ws.Range("AA9").Value = 1
    $ws.Range("AB9").Value = -1
    $ws.Range("AD9").Value = 0.4875
    $ws.Range("B288").Value = 8129704
    $ws.Range("E288").Value = "Botev Vratsa"
    $ws.Range("F288").Value = "Beroe"
    $ws.Range("G288").Value = 1
    $ws.Range("H288").Value = 0
    $ws.Range("I288").Value = 1
    $ws.Range("J288").Value = 0
    $ws.Range("L288").Value = 1.533
    $ws.Range("M288").Value = 3.9
    $ws.Range("N288").Value = 6
    $ws.Range("O288").Value = 1.4
    $ws.Range("P288").Value = 4.333
    $ws.Range("Q288").Value = 9
    $ws.Range("R288").Value = -1.25
    $ws.Range("S288").Value = 1.95
    $ws.Range("T288").Value = 1.9
    $ws.Range("V288").Value = 1.95
    $ws.Range("W288").Value = 1.9
    $ws.Range("X288").Value = 0.3999999999999999
    $ws.Range("AA288").Value = -0.5
    $ws.Range("AB288").Value = 0.45
    $ws.Range("AC288").Value = -1
    $ws.Range("AD288").Value = 0.8999999999999999
    $ws.Range("B289").Value = 8129703
    $ws.Range("E289").Value = "FC Hebar Pazardzhik"
    $ws.Range("F289").Value = "Etar 1924 Veliko Tarnovo"
    $ws.Range("G289").Value = 3
    $ws.Range("H289").Value = 1
    $ws.Range("J289").Value = 1
    $ws.Range("L289").Value = 1.333
    $ws.Range("M289").Value = 5
    $ws.Range("N289").Value = 8
    $ws.Range("O289").Value = 1.571
    $ws.Range("P289").Value = 4.2
    $ws.Range("Q289").Value = 5.25
    $ws.Range("R289").Value = -1
    $ws.Range("S289").Value = 2
    $ws.Range("T289").Value = 1.85
    $ws.Range("U289").Value = 2.5
    $ws.Range("V289").Value = 1.9
    $ws.Range("W289").Value = 1.95
    $ws.Range("X289").Value = 0.571
    $ws.Range("AA289").Value = 1
    $ws.Range("AC289").Value = 0.8999999999999999
    $ws.Range("AD289").Value = -1
    $ws.Range("B290").Value = 8129607
    $ws.Range("E290").Value = "Lokomotiv 1929 Sofia"
    $ws.Range("F290").Value = "Pirin Blagoevgrad"
    $ws.Range("I290").Value = 0
    $ws.Range("L290").Value = 5.75
    $ws.Range("M290").Value = 4
    $ws.Range("N290").Value = 1.533
    $ws.Range("O290").Value = 3.4
    $ws.Range("P290").Value = 3.5
    $ws.Range("Q290").Value = 2.05
    $ws.Range("R290").Value = 0.25
    $ws.Range("S290").Value = 2.025
    $ws.Range("T290").Value = 1.825
    $ws.Range("U290").Value = 2.25
    $ws.Range("V290").Value = 1.925
    $ws.Range("W290").Value = 1.925
    $ws.Range("X290").Value = 2.4
    $ws.Range("AA290").Value = 1.025
    $ws.Range("AB290").Value = -1
    $ws.Range("AD290").Value = 0.925
    $ws.Range("B294").Value = 8129620
    $ws.Range("E294").Value = "Slavia Sofia"
    $ws.Range("F294").Value = "CSKA 1948 Sofia"
    $ws.Range("H294").Value = 2
    $ws.Range("L294").Value = 3.75
    $ws.Range("M294").Value = 3.6
    $ws.Range("N294").Value = 1.9
    $ws.Range("O294").Value = 9.5
    $ws.Range("Q294").Value = 1.285
    $ws.Range("S294").Value = 1.925
    $ws.Range("T294").Value = 1.925
    $ws.Range("V294").Value = 1.875
    $ws.Range("W294").Value = 1.975
    $ws.Range("Z294").Value = 0.2849999999999999
    $ws.Range("AA294").Value = -1
    $ws.Range("AB294").Value = 0.925
    $ws.Range("AD294").Value = 0.9750000000000001
    $ws.Range("B295").Value = 8129621
    $ws.Range("E295").Value = "Botev Plovdiv"
    $ws.Range("F295").Value = "Arda Kardzhali"
    $ws.Range("H295").Value = 1
    $ws.Range("L295").Value = 4.333
    $ws.Range("M295").Value = 4
    $ws.Range("N295").Value = 1.727
    $ws.Range("O295").Value = 8.5
    $ws.Range("Q295").Value = 1.333
    $ws.Range("S295").Value = 1.875
    $ws.Range("T295").Value = 1.975
    $ws.Range("V295").Value = 1.85
    $ws.Range("W295").Value = 2
    $ws.Range("Z295").Value = 0.333
    $ws.Range("AA295").Value = 0.875
    $ws.Range("AB295").Value = -1
    $ws.Range("AD295").Value = 1
